$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "Hardware/software Kevin" header block from row 31 up to row 19,
# preserving its formatting (bold font + border on B/C/D).
$ws.Range("B31:D31").Copy($ws.Range("B19:D19"))
$ws.Range("B31:D31").Clear()

# Fill in Kevin's hardware / software list (rows 20-29)
$ws.Range("B20").Value = "Microsoft Windows 10 Home"
$ws.Range("D20").Value = "10.0.16299 Build 16299"

$ws.Range("B21").Value = "Intel(R) Core(TM) i5-4200U CPU @ 1.60GHz, 2301 MHz, 2 core('s), 4 logische processor(s)"
$ws.Range("C21").Value = "Processor"
$ws.Range("D21").Value = "n.v.t"

$ws.Range("B22").Value = "AMD Radeon R7 M265 Series"
$ws.Range("C22").Value = "Graphics Chipset"
$ws.Range("D22").Value = "n.v.t"

$ws.Range("B23").Value = "Wampserver"
$ws.Range("D23").Value = "3.0.6"

$ws.Range("B24").Value = "sublime"
$ws.Range("D24").Value = "build 3126"

$ws.Range("B25").Value = "nodeJS"
$ws.Range("D25").Value = "v8.6.0"

$ws.Range("B26").Value = "VeuJS"
$ws.Range("D26").Value = "2.5.16"

$ws.Range("B27").Value = "git"
$ws.Range("D27").Value = "2.10.0.windows.1"

$ws.Range("B28").Value = "gitkraken"
$ws.Range("D28").Value = "3.5.1"

$ws.Range("B29").Value = "Google Chrome"
$ws.Range("D29").Value = "65.0.3325.181"

# Also extend Tim's software list (rows 14-15) with git / gitkraken
$ws.Range("B14").Value = "git"
$ws.Range("B15").Value = "gitkraken"

$ws.Range("D28").Select()
